$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header, matching style of D1 (border + bold header style)
$ws.Range("E1").Value = "Pretrained"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# Values for column E (1 = pretrained, 0 = not pretrained)
$values = @(1,1,1,1,1,1,1,1,0,1,1,1,1,1,1,1,1,0,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

# Update the active selection to match the post-edit state (E21)
$ws.Range("E21").Select()
